# Rotate data for rows 7, 8, 9 in columns A,B,E,F,G,H,I,Q,R,AC:
#   new row 7 <- old row 9
#   new row 8 <- old row 7
#   new row 9 <- old row 8
# All other columns are identical across these rows and remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "I", "Q", "R", "AC")

# Capture current values for the rows involved before overwriting anything.
$old7 = @{}
$old8 = @{}
$old9 = @{}
foreach ($col in $cols) {
    $old7[$col] = $ws.Range("$col`7").Value()
    $old8[$col] = $ws.Range("$col`8").Value()
    $old9[$col] = $ws.Range("$col`9").Value()
}

foreach ($col in $cols) {
    $ws.Range("$col`7").Value = $old9[$col]
    $ws.Range("$col`8").Value = $old7[$col]
    $ws.Range("$col`9").Value = $old8[$col]
}
